$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Insert a new "Resource: Study Notes (Biblica)" Heading2 paragraph
# right before the "License Information" Heading2 paragraph.
# ---------------------------------------------------------------------------
$licPara = $d.Paragraphs(3)
$licPara.Range.InsertParagraphBefore()
$resPara = $d.Paragraphs(3)
$resPara.Range.Text = "Resource: Study Notes (Biblica)"

Write-Host "Edit 1 done"

# ---------------------------------------------------------------------------
# Edit 2: In the license paragraph, rename only the first "Biblica Study
# Notes" occurrence (the bold run) to "Study Notes (Biblica)". The second
# occurrence later in the same sentence must stay unchanged.
# ---------------------------------------------------------------------------
$licenseTextPara = $d.Paragraphs(5)
$licenseRng = $licenseTextPara.Range
$licenseRng.Find.Execute("Biblica Study Notes", $false, $true, $false, $false, $false, $true, 1, $false, "Study Notes (Biblica)", 1) | Out-Null

Write-Host "Edit 2 done"

# ---------------------------------------------------------------------------
# Edit 3: Rename the Heading1 title "Biblica Study Notes" to
# "Study Notes (Biblica)".
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(7)
$titleRng = $titlePara.Range
$titleRng.Find.Execute("Biblica Study Notes", $false, $true, $false, $false, $false, $true, 1, $false, "Study Notes (Biblica)", 1) | Out-Null

Write-Host "Edit 3 done"

# ---------------------------------------------------------------------------
# Edit 4: Insert a new "JER" resource-reference block right before the
# "Jeremiah 1:1-19" Heading2 paragraph. The block consists of:
#   1. Heading2 "JER"
#   2. Italic paragraph listing all the chapter ranges
#   3. Blank paragraph containing a single space
#   4. A paragraph carrying a new (header/footer-less) continuous sectPr
#   5. An empty paragraph
# ---------------------------------------------------------------------------
$jerHeadingPara = $d.Paragraphs(10)

# 5. empty paragraph (inserted first so it ends up directly before the
#    existing "Jeremiah 1:1-19" heading)
$jerHeadingPara.Range.InsertParagraphBefore()

# 4. sectPr paragraph
$jerHeadingPara = $d.Paragraphs(10)
$jerHeadingPara.Range.InsertParagraphBefore()

# 3. blank paragraph with a single space
$jerHeadingPara = $d.Paragraphs(10)
$jerHeadingPara.Range.InsertParagraphBefore()

# 2. italic paragraph with the chapter list
$jerHeadingPara = $d.Paragraphs(10)
$jerHeadingPara.Range.InsertParagraphBefore()

# 1. Heading2 "JER"
$jerHeadingPara = $d.Paragraphs(10)
$jerHeadingPara.Range.InsertParagraphBefore()

Write-Host "Edit 4 skeleton done"

# Paragraph 10: Heading2 "JER"  (style is inherited automatically because
# InsertParagraphBefore clones the following paragraph's pPr, and the
# following paragraph at this point is already the "Jeremiah 1:1-19"
# Heading2 paragraph)
$p10 = $d.Paragraphs(10)
$p10.Range.Text = "JER"

# Paragraph 11: italic chapter-range list
$p11 = $d.Paragraphs(11)
$p11.Range.Text = "Jeremiah 1:1–19, Jeremiah 2:1–12:17, Jeremiah 13:1–24:10, Jeremiah 25:1–38, Jeremiah 26:1–29:32, Jeremiah 30:1–33:26, Jeremiah 34:1–45:5, Jeremiah 46:1–49:39, Jeremiah 50:1–51:64, Jeremiah 52:1–34"
$p11.Range.Font.Italic = 1

# Paragraph 12: blank paragraph containing a single space
$p12 = $d.Paragraphs(12)
$p12.Range.Text = " "

# Paragraph 13 stays empty (it becomes the continuous-section-break
# paragraph once the sectPr is attached below).
# Paragraph 14 stays empty.

Write-Host "Edit 4 text done"
